# Daily attendance processing - 2025-12-29 08:42:33
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. "Recorded By" (column G) values: swap "System, <email>" -> "<email>, System"
# Column G = 7
$recordedByRows = @(8,9,10,12,14,15,17,34,35,36,38,40,41,43,60,61,62,64,66,67,69,86,87,88,90,92,93,95,112,113,114,116,118,119,121,138,139,140,142,144,145,147,164,167,170,191,194,197,218,221,224,245,248,251,272,275,278,299,302,305)
foreach ($r in $recordedByRows) {
    $ws.Cells.Item($r, 7).Value = "dnasr281@gmail.com, System"
}

# --- 2. Newly-processed sessions: rows 18, 44, 70, 96, 122, 148
#     status flips from "Pending" (yellow fill) to "Not Recorded" (its own fill).
#     Row 3 already carries the target "Not Recorded" style, so copy its
#     format (A3:I3) onto the target row's A:I cells, then fix up the Status text.
$notRecordedRows = @(18,44,70,96,122,148)
$ws.Range("A3:I3").Copy() | Out-Null
foreach ($r in $notRecordedRows) {
    $targetRange = "A" + $r + ":I" + $r
    $ws.Range($targetRange).PasteSpecial(-4122) | Out-Null
    $ws.Cells.Item($r, 9).Value = "Not Recorded"
}
$excel.CutCopyMode = 0

# --- 3. Class Statistics block (K:L columns) - Missing/Pending session counters
$ws.Range("L7").Value = 9
$ws.Range("L8").Value = 126

# --- 4. Per-group summary block (columns K:S) - Late (P) / OnTime-ish (Q) counters
$ws.Range("P15").Value = 2
$ws.Range("Q15").Value = 9

$ws.Range("P16").Value = 1
$ws.Range("Q16").Value = 9

$ws.Range("P17").Value = 1
$ws.Range("Q17").Value = 9

$ws.Range("P18").Value = 1
$ws.Range("Q18").Value = 9

$ws.Range("P19").Value = 1
$ws.Range("Q19").Value = 9

$ws.Range("P20").Value = 2
$ws.Range("Q20").Value = 9
